# Adds four new daily exchange-rate rows (13, 14, 15 and 18 October 2021)
# to the bottom of the "Tipos de cambio 2021 - Diaria" sheet, matching the
# "Actualizacion desde MV -datos-" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the 60 values (columns A..BH) for one new row,
# in the same left-to-right column order as the existing data rows.
$newRows = @(
    @("13-10-2021", 33.306, 1, 4.1537, 6.86, 628.55, 22.0315, 6.4555, 129.97, 8.5602, 8.7826, 0.7097, 3.672, 9.08, 1.3604, 1.2465, 1, 2.1017, 0.82, 1, 1.3572, 7.7813, 1.4424, 28.155, 22734, 0.8677, 312.57, 102.8307, 0.9312, 6880, 26.1143, 4.2932, 15.6935, 0.7363, 9.034700000000001, 4.0555, 99.09099999999999, 827.5599999999999, 3727.32, 24, 56.045, 50.835, 20.7613, 43.4046, 7.7183, 14.9639, 5.5373, 42000, 3.75, 4.163, 71.8329, 75.51000000000001, 14215, 170.9846, 3.2301, 424.64, 1198.72, 113.65, 6.4558, 3.9732),
    @("14-10-2021", 33.173, 1, 4.1528, 6.86, 628.48, 21.9075, 6.4212, 128.84, 8.4953, 8.6957, 0.7095, 3.672, 9.0291, 1.3552, 1.2435, 1, 2.0995, 0.82, 1, 1.3522, 7.7793, 1.4364, 28.092, 22736, 0.863, 310.67, 102.711, 0.9243, 6890, 26.3399, 4.2695, 15.6537, 0.7323, 9.066000000000001, 3.9872, 99.1211, 820.1900000000001, 3734.8, 24, 56.06, 50.655, 20.6412, 43.62, 7.7053, 14.8014, 5.57, 42000, 3.75, 4.159, 71.97969999999999, 75.37, 14215, 171, 3.2294, 425.57, 1193.56, 113.36, 6.4256, 3.9508),
    @("15-10-2021", 33.169, 1, 4.1711, 6.86, 628.0700000000001, 21.8998, 6.4201, 128.98, 8.4617, 8.6424, 0.708, 3.672, 9.053800000000001, 1.3495, 1.2372, 1, 2.0964, 0.82, 1, 1.3491, 7.7783, 1.4239, 28.089, 22736, 0.8629, 309.36, 102.4964, 0.9237, 6895, 26.3399, 4.2701, 15.6537, 0.731, 9.1843, 3.9207, 99.1508, 816.8099999999999, 3765.34, 24, 56.24, 50.6, 20.5732, 43.735, 7.7271, 14.795, 5.5111, 42000, 3.75, 4.1525, 71.4169, 75.2525, 14115, 171, 3.221, 425.46, 1186.52, 113.68, 6.4368, 3.941),
    @("18-10-2021", 33.396, 1, 4.1515, 6.86, 625.28, 21.868, 6.4158, 128.71, 8.4245, 8.610300000000001, 0.7083, 3.672, 9.042999999999999, 1.3484, 1.2385, 1, 2.0886, 0.82, 1, 1.3485, 7.7791, 1.4144, 27.985, 22734, 0.8623, 309.76, 102.5114, 0.9236, 6905, 26.3433, 4.2638, 15.6845, 0.727, 9.222899999999999, 3.9305, 99.18000000000001, 826.1900000000001, 3766.5, 24, 56.26, 50.708, 20.3738, 43.79, 7.7031, 14.6227, 5.4427, 42000, 3.75, 4.156, 71.01130000000001, 75.2525, 14070, 171.0862, 3.2199, 425.61, 1182.2, 114.27, 6.4343, 3.934)
)

$startRow = 197
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
